# Add a "Save" column (H) to the s_vals sheet.
#
# Column H is a 0/1 indicator derived from column G ("sum"): it is 1 when
# the combined "sum" value for that outing reached the high tier the
# underlying data uses to flag a save (>= ~8.4186), and 0 otherwise.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 66

# Header cell, matching the style used by the other header cells (B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Threshold that separates "save" outings from the rest, based on the
# "sum" (column G) value of the outing.
$threshold = 8.4

$count = $lastRow - 1
$values = New-Object 'object[,]' $count,1

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -ge $threshold) {
        $save = 1
    } else {
        $save = 0
    }
    $values[$r - 2, 0] = $save
}

$ws.Range("H2:H$lastRow").Value = $values

$ws.Range("A1:H$lastRow").Select()
